$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FINAL_STATESUMMARY")

$ws.Range("A76").Copy()
$ws.Range("A77").PasteSpecial(-4122)

$ws.Range("A77").Value = 44257
$ws.Range("B77").Value = 1435667
$ws.Range("C77").Value = 42333
$ws.Range("D77").Value = 933250
$ws.Range("E77").Value = 27105
$ws.Range("F77").Value = 498245
$ws.Range("G77").Value = 15187

$wb.Names.Item("FINAL_STATESUMMARY").RefersTo = "='FINAL_STATESUMMARY'!`$A`$1:`$G`$77"
